$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3399.1714
$ws.Range("I19").Value = 6454.0625
$ws.Range("J19").Value = 826.6316
$ws.Range("K19").Value = 6454.0625
$ws.Range("L19").Value = 826.6316
$ws.Range("M19").Value = -6279.0625
$ws.Range("N19").Value = -1176.6316

$ws.Range("H132").Value = 2022.3721
$ws.Range("I132").Value = 1671.1666
$ws.Range("J132").Value = 3828.5715
$ws.Range("K132").Value = 5013.4998
$ws.Range("L132").Value = 11485.7145
$ws.Range("M132").Value = -2483.4998
$ws.Range("N132").Value = -16545.7145

$ws.Range("H137").Value = 2034.4036
$ws.Range("I137").Value = 1885.3954
$ws.Range("J137").Value = 2492.0715
$ws.Range("K137").Value = 5656.1862
$ws.Range("L137").Value = 7476.2145
$ws.Range("M137").Value = -3106.1862
$ws.Range("N137").Value = -12576.2145

$ws.Range("H138").Value = 2073.7402
$ws.Range("I138").Value = 936.6667
$ws.Range("J138").Value = 3072.1462
$ws.Range("K138").Value = 2810.0001
$ws.Range("L138").Value = 9216.438600000001
$ws.Range("M138").Value = 2329.9999
$ws.Range("N138").Value = -19496.4386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1128782
$ws.Range("I32").Value = 1284169.1
$ws.Range("J32").Value = 6542
$ws.Range("K32").Value = 1284169.1
$ws.Range("L32").Value = 6542
$ws.Range("M32").Value = -1283882.1
$ws.Range("N32").Value = -7116

$ws.Range("H74").Value = 1608.8853
$ws.Range("I74").Value = 1012.53064
$ws.Range("J74").Value = 4044
$ws.Range("K74").Value = 1012.53064
$ws.Range("L74").Value = 4044
$ws.Range("M74").Value = -138.5306399999999
$ws.Range("N74").Value = -5792

$ws.Range("H77").Value = 1608.8853
$ws.Range("I77").Value = 1012.53064
$ws.Range("J77").Value = 4044
$ws.Range("K77").Value = 5062.6532
$ws.Range("L77").Value = 20220
$ws.Range("M77").Value = -694.6531999999997
$ws.Range("N77").Value = -28956

$ws.Range("H110").Value = 1236.5
$ws.Range("I110").Value = 1620.4445
$ws.Range("J110").Value = 742.8570999999999
$ws.Range("K110").Value = 1620.4445
$ws.Range("L110").Value = 742.8570999999999
$ws.Range("M110").Value = 424.5554999999999
$ws.Range("N110").Value = -4832.8571

$ws.Range("H132").Value = 23297.348
$ws.Range("I132").Value = 27464.3
$ws.Range("J132").Value = 4777.5557
$ws.Range("K132").Value = 82392.89999999999
$ws.Range("L132").Value = 14332.6671
$ws.Range("M132").Value = -79862.89999999999
$ws.Range("N132").Value = -19392.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1904.4667
$ws.Range("I20").Value = 1734.1666
$ws.Range("J20").Value = 2018
$ws.Range("K20").Value = 1734.1666
$ws.Range("L20").Value = 2018
$ws.Range("M20").Value = -1487.1666
$ws.Range("N20").Value = -2512

$ws.Range("H86").Value = 7317.7
$ws.Range("I86").Value = 17866.666
$ws.Range("J86").Value = 2796.7144
$ws.Range("K86").Value = 17866.666
$ws.Range("L86").Value = 2796.7144
$ws.Range("M86").Value = -16743.666
$ws.Range("N86").Value = -5042.7144

$ws.Range("H89").Value = 7317.7
$ws.Range("I89").Value = 17866.666
$ws.Range("J89").Value = 2796.7144
$ws.Range("K89").Value = 89333.33
$ws.Range("L89").Value = 13983.572
$ws.Range("M89").Value = -83717.33
$ws.Range("N89").Value = -25215.572

$ws.Range("H134").Value = 1654.807
$ws.Range("I134").Value = 1170.7234
$ws.Range("J134").Value = 3930
$ws.Range("K134").Value = 3512.1702
$ws.Range("L134").Value = 11790
$ws.Range("M134").Value = -977.1702000000005
$ws.Range("N134").Value = -16860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2250.8572
$ws.Range("I31").Value = 1551.32
$ws.Range("J31").Value = 3999.7
$ws.Range("K31").Value = 1551.32
$ws.Range("L31").Value = 3999.7
$ws.Range("M31").Value = -1256.32
$ws.Range("N31").Value = -4589.7

$ws.Range("H34").Value = 2250.8572
$ws.Range("I34").Value = 1551.32
$ws.Range("J34").Value = 3999.7
$ws.Range("K34").Value = 1551.32
$ws.Range("L34").Value = 3999.7
$ws.Range("M34").Value = -1349.32
$ws.Range("N34").Value = -4403.7

$ws.Range("H58").Value = 1152.4259
$ws.Range("I58").Value = 641.2381
$ws.Range("J58").Value = 2941.5833
$ws.Range("K58").Value = 641.2381
$ws.Range("L58").Value = 2941.5833
$ws.Range("M58").Value = -438.2381
$ws.Range("N58").Value = -3347.5833

$ws.Range("H134").Value = 1540.5676
$ws.Range("I134").Value = 922.1852
$ws.Range("J134").Value = 3210.2
$ws.Range("K134").Value = 2766.5556
$ws.Range("L134").Value = 9630.599999999999
$ws.Range("M134").Value = -231.5556000000001
$ws.Range("N134").Value = -14700.6

$ws.Range("H136").Value = 1152.4259
$ws.Range("I136").Value = 641.2381
$ws.Range("J136").Value = 2941.5833
$ws.Range("K136").Value = 1923.7143
$ws.Range("L136").Value = 8824.749899999999
$ws.Range("M136").Value = 626.2856999999999
$ws.Range("N136").Value = -13924.7499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 11875
$ws.Range("J39").Value = 11875
$ws.Range("L39").Value = 11875
$ws.Range("N39").Value = -12939

$ws.Range("H70").Value = 6875.2764
$ws.Range("I70").Value = 3960.7144
$ws.Range("J70").Value = 9229.346
$ws.Range("K70").Value = 3960.7144
$ws.Range("L70").Value = 9229.346
$ws.Range("M70").Value = -3690.7144
$ws.Range("N70").Value = -9769.346

$ws.Range("H73").Value = 6875.2764
$ws.Range("I73").Value = 3960.7144
$ws.Range("J73").Value = 9229.346
$ws.Range("K73").Value = 3960.7144
$ws.Range("L73").Value = 9229.346
$ws.Range("M73").Value = -3024.7144
$ws.Range("N73").Value = -11101.346

$ws.Range("H122").Value = 1428.1852
$ws.Range("I122").Value = 1361.7222
$ws.Range("J122").Value = 1561.1111
$ws.Range("K122").Value = 4085.1666
$ws.Range("L122").Value = 4683.3333
$ws.Range("M122").Value = -1635.1666
$ws.Range("N122").Value = -9583.3333

$ws.Range("H132").Value = 3077
$ws.Range("I132").Value = 2798.279
$ws.Range("J132").Value = 4075.75
$ws.Range("K132").Value = 8394.837
$ws.Range("L132").Value = 12227.25
$ws.Range("M132").Value = -5864.837
$ws.Range("N132").Value = -17287.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1548.75
$ws.Range("I61").Value = 895
$ws.Range("J61").Value = 1766.6666
$ws.Range("K61").Value = 895
$ws.Range("L61").Value = 1766.6666
$ws.Range("M61").Value = -693
$ws.Range("N61").Value = -2170.6666

$ws.Range("H100").Value = 62504240
$ws.Range("I100").Value = 4696
$ws.Range("K100").Value = 4696
$ws.Range("M100").Value = -4155

$ws.Range("H113").Value = 1548.75
$ws.Range("I113").Value = 895
$ws.Range("J113").Value = 1766.6666
$ws.Range("K113").Value = 895
$ws.Range("L113").Value = 1766.6666
$ws.Range("M113").Value = 1275
$ws.Range("N113").Value = -6106.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 534.46155
$ws.Range("I113").Value = 636.2222
$ws.Range("J113").Value = 305.5
$ws.Range("K113").Value = 1908.6666
$ws.Range("L113").Value = 916.5
$ws.Range("M113").Value = 261.3334
$ws.Range("N113").Value = -5256.5

$ws.Range("H126").Value = 1059.3667
$ws.Range("J126").Value = 1491.25
$ws.Range("L126").Value = 4473.75
$ws.Range("N126").Value = -9413.75

$ws.Range("H132").Value = 1705.7838
$ws.Range("I132").Value = 1203.742
$ws.Range("J132").Value = 4299.6665
$ws.Range("K132").Value = 3611.226
$ws.Range("L132").Value = 12898.9995
$ws.Range("M132").Value = -1081.226
$ws.Range("N132").Value = -17958.9995

$ws.Range("H136").Value = 26318596
$ws.Range("I136").Value = 38463892
$ws.Range("J136").Value = 3791.6667
$ws.Range("K136").Value = 115391676
$ws.Range("L136").Value = 11375.0001
$ws.Range("M136").Value = -115389126
$ws.Range("N136").Value = -16475.0001

$ws.Range("H138").Value = 35824.75
$ws.Range("J138").Value = 35824.75
$ws.Range("L138").Value = 35824.75
$ws.Range("N138").Value = -46104.75

$ws.Range("H141").Value = 30285.77
$ws.Range("J141").Value = 30285.77
$ws.Range("L141").Value = 30285.77
$ws.Range("N141").Value = -40645.77
